# product-management.xlsx — "6. Đẩy file tĩnh lên trên cloud" row added
# (commit message: "1. Sắp xếp sản phẩm theo các tiêu chí khác nhau")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row 21: "6. Đẩy file tĩnh lên trên cloud" + commit link ---
$ws.Range("B21").Value = "6. Đẩy file tĩnh lên trên cloud"
$ws.Range("C21").Value = "https://github.com/nguyentienminh07102004/product-management/commit/40564ed629ff1085a1e16ddb2c9fdef618240676"

# Turn C21 into a real hyperlink, same as the other "commit link" cells above it.
$ws.Hyperlinks.Add($ws.Cells.Item(21, 3), "https://github.com/nguyentienminh07102004/product-management/commit/40564ed629ff1085a1e16ddb2c9fdef618240676") | Out-Null

# Match the visual style Excel uses for the other link cells (C19/C20) exactly,
# instead of the composite style Hyperlinks.Add would otherwise synthesize.
$ws.Range("C19").Copy()
$ws.Range("C21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- View state: scroll/select like the saved workbook (column C in view, C23 selected) ---
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("C23").Select()
